# Update cryptos list - price/volume refresh as of Mon May 13 09:26:10 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.936.71"
$ws.Range("E2").Value = "  +3.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.977.92"
$ws.Range("E3").Value = "  +2.13%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'596.79"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'145.63"
$ws.Range("E6").Value = "  +0.47%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "2.975.50"
$ws.Range("E8").Value = "  +2.07%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.05%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'7.38"
$ws.Range("E10").Value = "  +6.31%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = "  +2.47%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = "  +2.75%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "  +4.62%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'33.53"
$ws.Range("E14").Value = "  +0.46%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.469.69"
$ws.Range("E16").Value = "  +2.13%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.818.64"
$ws.Range("E17").Value = "  +3.18%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "'6.75"
$ws.Range("E18").Value = "  +1.01%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.948.46"
$ws.Range("E19").Value = "  +1.15%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'444.05"
$ws.Range("E20").Value = "  +2.56%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'13.58"
$ws.Range("E21").Value = "  +1.65%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.678"
$ws.Range("E22").Value = "  +0.51%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'7.17"
$ws.Range("E23").Value = "  +0.83%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'82.22"
$ws.Range("E24").Value = "  +1.00%  "

# Row 25 - RenderToken
$ws.Range("D25").Value = "'10.87"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'12.09"
$ws.Range("E26").Value = "  +2.75%  "

# Row 27 - Fetch.AI
$ws.Range("D27").Value = "'2.15"
$ws.Range("E27").Value = "  -2.47%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.00%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +1.33%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "'7.08"
$ws.Range("E30").Value = "  +1.81%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -6.16%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'26.57"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.28%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  +0.01%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0₃0882"
$ws.Range("E35").Value = "  +1.72%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "  -1.84%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  +0.84%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  -1.63%  "

# Row 41 - Cosmos
$ws.Range("D41").Value = "'8.64"
$ws.Range("E41").Value = "  +1.05%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -2.38%  "

# Row 43 - TheGraph
$ws.Range("D43").Value = "'0.287"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44 - Arweave
$ws.Range("D44").Value = "'39.01"
$ws.Range("E44").Value = "  -5.97%  "

# Row 45 - Bittensor
$ws.Range("D45").Value = "'371.72"
$ws.Range("E45").Value = "  -1.02%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.706.71"
$ws.Range("E46").Value = "  +0.59%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  -0.76%  "

# Row 48 - Monero
$ws.Range("E48").Value = "  +0.88%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'23.39"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.24%  "

# Row 38 - now Stacks (was OKB)
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.05"
$ws.Range("E38").Value = "  +4.01%  "

# Row 39 - now OKB (was Stacks)
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'49.64"
$ws.Range("E39").Value = "  +0.20%  "
